$d = $word.ActiveDocument

$rng = $d.Content
$rng.Find.Execute("  --count 3 \", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if ($rng.Find.Found) {
    $full = $rng.Text
    $prefixLen = $full.IndexOf("3")
    $start = $rng.Start
    $midStart = $start + $prefixLen
    $midEnd = $midStart + 1

    $midRange = $d.Range($midStart, $midEnd)
    $midRange.Text = "2"
}
